$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Match the header style used by the existing header cells (e.g. H1),
# then set the header text for the new I0 / IF columns.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for I2:J82
$data = @(
    @(8, 8),
    @(8, 8),
    @(7, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(7, 8),
    @(9, 9),
    @(8, 8),
    @(9, 9),
    @(8, 8),
    @(8, 8),
    @(11, 11),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(9, 9),
    @(8, 8),
    @(8, 8),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(8, 8),
    @(8, 8),
    @(9, 9),
    @(8, 8),
    @(9, 9),
    @(9, 9),
    @(7, 7),
    @(7, 8),
    @(9, 9),
    @(9, 9),
    @(8, 8),
    @(7, 8),
    @(8, 8),
    @(9, 9),
    @(7, 8),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(9, 9),
    @(8, 8),
    @(8, 8),
    @(9, 9),
    @(8, 8),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(8, 9),
    @(8, 8),
    @(9, 9),
    @(8, 8),
    @(7, 7),
    @(7, 7),
    @(8, 8),
    @(9, 9),
    @(7, 7),
    @(8, 8),
    @(8, 8),
    @(6, 6),
    @(8, 8),
    @(9, 9),
    @(9, 9),
    @(8, 8),
    @(9, 9),
    @(9, 9),
    @(7, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(7, 7),
    @(5, 5),
    @(7, 7)
)

for ($idx = 0; $idx -lt $data.Count; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $data[$idx][0]
    $ws.Cells.Item($row, 10).Value = $data[$idx][1]
}
